$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / link / percentage cells (never numeric-looking, safe as literal text)
$ws.Range("D2").Value = '26.717.95'
$ws.Range("E2").Value = '  -0.21%  '
$ws.Range("D3").Value = '1.638.54'
$ws.Range("E3").Value = '  -0.60%  '
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("E5").Value = '  +0.66%  '
$ws.Range("E6").Value = '  -0.86%  '
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("E8").Value = '  -0.50%  '
$ws.Range("E9").Value = '  -0.71%  '
$ws.Range("E10").Value = '  -0.41%  '
$ws.Range("E11").Value = '  +0.38%  '
$ws.Range("D12").Value = '1.865.48'
$ws.Range("E12").Value = '  -0.74%  '
$ws.Range("D13").Value = '1.622.95'
$ws.Range("E13").Value = '  -1.61%  '
$ws.Range("E14").Value = '  -1.35%  '
$ws.Range("E15").Value = '  -1.44%  '
$ws.Range("E16").Value = '  -1.57%  '
$ws.Range("D17").Value = '26.704.98'
$ws.Range("E17").Value = '  -0.35%  '
$ws.Range("E18").Value = '  -2.47%  '
$ws.Range("E19").Value = '  -3.31%  '
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("E21").Value = '  -0.76%  '
$ws.Range("E22").Value = '  -1.22%  '
$ws.Range("E23").Value = '  -3.43%  '
$ws.Range("E24").Value = '  -2.72%  '
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("E26").Value = '  -0.14%  '
$ws.Range("E27").Value = '  -1.97%  '
$ws.Range("E28").Value = '  -0.39%  '
$ws.Range("E29").Value = '  -1.14%  '
$ws.Range("E30").Value = '  -2.88%  '
$ws.Range("E31").Value = '  +0.39%  '
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("E33").Value = '  -1.09%  '
$ws.Range("D34").Value = '1.269.56'
$ws.Range("E34").Value = '  -1.21%  '
$ws.Range("E35").Value = '  -1.11%  '
$ws.Range("E36").Value = '  -0.60%  '
$ws.Range("E37").Value = '  -2.32%  '
$ws.Range("E38").Value = '  -1.79%  '
$ws.Range("E39").Value = '  -2.69%  '
$ws.Range("E40").Value = '  -0.06%  '
$ws.Range("E41").Value = '  -1.26%  '
$ws.Range("E42").Value = '  -2.96%  '
$ws.Range("E43").Value = '  -3.71%  '
$ws.Range("D44").Value = '1.775.96'
$ws.Range("E44").Value = '  -0.83%  '
$ws.Range("E45").Value = '  -0.83%  '
$ws.Range("E46").Value = '  +0.93%  '
$ws.Range("E47").Value = '  -1.66%  '
$ws.Range("E48").Value = '  +0.34%  '
$ws.Range("E49").Value = '  -3.41%  '
$ws.Range("E50").Value = '  -1.15%  '
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("E51").Value = '  -0.13%  '

# Price cells whose new values look like plain numbers (single decimal point).
# Excel would silently coerce these to numeric cells (losing the exact text, e.g.
# trailing zeros / "0.0960" -> 9.6E-02), so force the Text number format first,
# then drop back to the unstyled "Normal" cell style so no stray style lingers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.80'
$ws.Range("D5").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0623'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.09'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0845'
$ws.Range("D11").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.49'
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '211.47'
$ws.Range("D19").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.32'
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.74'
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.09'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.56'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0503'
$ws.Range("D30").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.529'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.806'
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.803'
$ws.Range("D41").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '91.29'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.25'
$ws.Range("D46").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0960'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.01'
$ws.Range("D51").Style = "Normal"
